# Refresh the cryptocurrency price/volume snapshot in columns D (Price) and E (Volume 1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.434.98"
$ws.Range("D3").Value = "1.656.27"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'213.28"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'0.539"
$ws.Range("E6").Value = "  +5.09%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'23.50"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "'0.261"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("D12").Value = "1.889.43"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.660.64"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("D16").Value = "'65.47"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "27.434.01"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "'229.19"
$ws.Range("E18").Value = "  -7.07%  "
$ws.Range("D19").Value = "0.0₃0727"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").Value = "'9.38"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").Value = "'146.93"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'7.07"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'15.62"
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("D30").Value = "'0.0493"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("E31").Value = "  -4.36%  "
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").Value = "'3.13"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "1.423.96"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'0.906"
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "'1.04"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").Value = "'65.25"
$ws.Range("E43").Value = "  -5.61%  "
$ws.Range("D44").Value = "'2.22"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "'0.792"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "1.798.71"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("D48").Value = "'87.94"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "'7.71"
$ws.Range("E51").Value = "  -1.43%  "
